# Scheduled runner: refresh market-board derived figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) in columns H:N across all job sheets.
$wb = $excel.ActiveWorkbook

# ============ Sheet ALC ============
$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 521.76666
$ws.Range("I19").Value = 432.52942
$ws.Range("J19").Value = 638.46155
$ws.Range("K19").Value = 432.52942
$ws.Range("L19").Value = 638.46155
$ws.Range("M19").Value = -257.52942
$ws.Range("N19").Value = -988.46155
# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 55556656
$ws.Range("I40").Value = 166667170
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 166667170
$ws.Range("L40").Value = 1400
$ws.Range("M40").Value = -166666995
$ws.Range("N40").Value = -1750
# Row 55: A Real Smooth Move / Lanolin
$ws.Range("H55").Value = 184.4
$ws.Range("I55").Value = 100
$ws.Range("K55").Value = 100
$ws.Range("M55").Value = 114
# Row 61: Not Taking No for an Answer / Mega-Potion of Strength
$ws.Range("H61").Value = 98
$ws.Range("I61").Value = 98
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 294
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -122
$ws.Range("N61").ClearContents()
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 404121.53
$ws.Range("I132").Value = 459084.47
$ws.Range("J132").Value = 1060
$ws.Range("K132").Value = 1377253.41
$ws.Range("L132").Value = 3180
$ws.Range("M132").Value = -1374723.41
$ws.Range("N132").Value = -8240
# Row 133: Big Brush, Big Dreams / Ginseng Angle Brush
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120

# ============ Sheet ARM ============
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 1731799.9
$ws.Range("I2").Value = 1943.3334
$ws.Range("J2").Value = 3677888.5
$ws.Range("K2").Value = 1943.3334
$ws.Range("L2").Value = 3677888.5
$ws.Range("M2").Value = -1830.3334
$ws.Range("N2").Value = -3678114.5
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 9936.571
$ws.Range("I32").Value = 7776.884
$ws.Range("J32").Value = 17080.154
$ws.Range("K32").Value = 7776.884
$ws.Range("L32").Value = 17080.154
$ws.Range("M32").Value = -7489.884
$ws.Range("N32").Value = -17654.154
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 1742.174
$ws.Range("I61").Value = 1110.5897
$ws.Range("J61").Value = 5261
$ws.Range("K61").Value = 1110.5897
$ws.Range("L61").Value = 5261
$ws.Range("M61").Value = -898.5897
$ws.Range("N61").Value = -5685
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 3041.38
$ws.Range("I74").Value = 503.75
$ws.Range("K74").Value = 503.75
$ws.Range("M74").Value = 370.25
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 3041.38
$ws.Range("I77").Value = 503.75
$ws.Range("K77").Value = 2518.75
$ws.Range("M77").Value = 1849.25
# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 919.14703
$ws.Range("I110").Value = 811.3043
$ws.Range("K110").Value = 811.3043
$ws.Range("M110").Value = 1233.6957
# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 1731799.9
$ws.Range("I116").Value = 1943.3334
$ws.Range("J116").Value = 3677888.5
$ws.Range("K116").Value = 1943.3334
$ws.Range("L116").Value = 3677888.5
$ws.Range("M116").Value = 350.6666
$ws.Range("N116").Value = -3682476.5
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1600.5209
$ws.Range("I132").Value = 1472.2162
$ws.Range("K132").Value = 4416.6486
$ws.Range("M132").Value = -1886.6486
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1742.174
$ws.Range("I136").Value = 1110.5897
$ws.Range("J136").Value = 5261
$ws.Range("K136").Value = 3331.7691
$ws.Range("L136").Value = 15783
$ws.Range("M136").Value = -781.7691
$ws.Range("N136").Value = -20883

# ============ Sheet BSM ============
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 1731799.9
$ws.Range("I3").Value = 1943.3334
$ws.Range("J3").Value = 3677888.5
$ws.Range("K3").Value = 1943.3334
$ws.Range("L3").Value = 3677888.5
$ws.Range("M3").Value = -1829.3334
$ws.Range("N3").Value = -3678116.5
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 2253.25
$ws.Range("I20").Value = 2121.4707
$ws.Range("K20").Value = 2121.4707
$ws.Range("M20").Value = -1874.4707

# ============ Sheet CRP ============
$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 53.117645
$ws.Range("I7").Value = 49.3
$ws.Range("J7").Value = 58.57143
$ws.Range("K7").Value = 49.3
$ws.Range("L7").Value = 58.57143
$ws.Range("M7").Value = 63.7
$ws.Range("N7").Value = -284.57143
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1323.3334
$ws.Range("I31").Value = 1244.7826
$ws.Range("K31").Value = 1244.7826
$ws.Range("M31").Value = -949.7826
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1323.3334
$ws.Range("I34").Value = 1244.7826
$ws.Range("K34").Value = 1244.7826
$ws.Range("M34").Value = -1042.7826
# Row 63: So You Think You Can Lance? / Mythrite Trident
$ws.Range("H63").Value = 63600
$ws.Range("J63").Value = 63600
$ws.Range("L63").Value = 63600
$ws.Range("N63").Value = -64972
# Row 66: Sticks and Stones (L) / Mythrite Trident
$ws.Range("H66").Value = 63600
$ws.Range("J66").Value = 63600
$ws.Range("L66").Value = 190800
$ws.Range("N66").Value = -197664
# Row 92: Walk the Walk / Beech Rod
$ws.Range("H92").Value = 47500
$ws.Range("J92").Value = 47500
$ws.Range("L92").Value = 47500
$ws.Range("N92").Value = -52492

# ============ Sheet CUL ============
$ws = $wb.Worksheets.Item("CUL")
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1925781.6
$ws.Range("J131").Value = 2502155.5
$ws.Range("L131").Value = 7506466.5
$ws.Range("N131").Value = -7516546.5
# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 50000870
$ws.Range("I132").Value = 52632380
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 473691420
$ws.Range("L132").Value = 19800
$ws.Range("M132").Value = -473688890
$ws.Range("N132").Value = -24860

# ============ Sheet GSM ============
$ws = $wb.Worksheets.Item("GSM")
# Row 95: Chain of Command / Koppranickel Temple Chain
$ws.Range("H95").Value = 42672
$ws.Range("J95").Value = 42672
$ws.Range("L95").Value = 42672
$ws.Range("N95").Value = -48164
# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 726.2
$ws.Range("I107").Value = 428.5
$ws.Range("J107").Value = 834.4545000000001
$ws.Range("K107").Value = 428.5
$ws.Range("L107").Value = 834.4545000000001
$ws.Range("M107").Value = 1491.5
$ws.Range("N107").Value = -4674.4545

# ============ Sheet LTW ============
$ws = $wb.Worksheets.Item("LTW")
# Row 57: Too Hot to Handle / Raptorskin Wristbands
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 15200.1
$ws.Range("I122").Value = 28753
$ws.Range("J122").Value = 6164.8335
$ws.Range("K122").Value = 86259
$ws.Range("L122").Value = 18494.5005
$ws.Range("M122").Value = -83809
$ws.Range("N122").Value = -23394.5005

# ============ Sheet WVR ============
$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1827.4
$ws.Range("I122").Value = 1793.5
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 5380.5
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -2930.5
$ws.Range("N122").Value = -10450
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 3344.2222
$ws.Range("J132").Value = 4019.6
$ws.Range("L132").Value = 12058.8
$ws.Range("N132").Value = -17118.8
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 15761.857
$ws.Range("I136").Value = 18005.916
$ws.Range("K136").Value = 54017.74800000001
$ws.Range("M136").Value = -51467.74800000001
